$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2, shifting existing data down.
$ws.Rows.Item(2).Insert()

# Fill in the new test line in the newly inserted row 2.
# Note: shared-string table entries are appended in the order values are
# first assigned, so E2 is set before B2 to reproduce the original
# uniqueCount ordering (39 = "ad dies..." / 40 = "...tertius.").
$ws.Cells.Item(2, 1).Value = 1200000008
$ws.Cells.Item(2, 5).Value = "ad dies ist ein längeres Zitat, das so auch im anderen Text auftaucht und von Tracer gefunden werden müsste."
$ws.Cells.Item(2, 2).Value = "dies ist ein längeres Zitat, das so auch im anderen Text auftaucht und von Tracer gefunden werden müsste tertius."
$ws.Cells.Item(2, 4).Value = 1100000013

# Match the saved selection state from the authored workbook.
$ws.Range("B3").Select()
